$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '67.638.75'
$ws.Range("E2").Value = '  -1.58%  '
$ws.Range("D3").Value = '2.427.43'
$ws.Range("E3").Value = '  -1.16%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'552.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("D6").Value = "'159.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.15%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  +1.01%  '
$ws.Range("E9").Value = '  +6.35%  '
$ws.Range("E10").Value = '  -0.70%  '
$ws.Range("E11").Value = '  -1.30%  '
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").Value = '67.585.52'
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("E14").Value = '  +1.14%  '
$ws.Range("D15").Value = "'23.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.21%  '
$ws.Range("E16").Value = '  -3.25%  '
$ws.Range("D17").Value = "'329.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.29%  '
$ws.Range("E18").Value = '  -2.11%  '
$ws.Range("D19").Value = "'3.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("D22").Value = "'66.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.92%  '
$ws.Range("E23").Value = '  -1.07%  '
$ws.Range("D24").Value = "'8.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.36%  '
$ws.Range("D25").Value = '0.0₃0804'
$ws.Range("E25").Value = '  -1.35%  '
$ws.Range("E26").Value = '  -1.62%  '
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("D28").Value = "'416.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.89%  '
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("D31").Value = "'160.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.01%  '
$ws.Range("E32").Value = '  -0.59%  '
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("D34").Value = "'17.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.22%  '
$ws.Range("E35").Value = '  -3.48%  '
$ws.Range("D36").Value = "'0.295"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.34%  '
$ws.Range("E37").Value = '  -3.54%  '
$ws.Range("D38").Value = "'1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.74%  '
$ws.Range("E39").Value = '  -2.56%  '
$ws.Range("E40").Value = '  -2.98%  '
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").Value = "'129.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("D43").Value = "'0.0707"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.69%  '
$ws.Range("D44").Value = "'0.476"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.95%  '
$ws.Range("D45").Value = "'0.554"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.22%  '
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("E48").Value = '  -7.11%  '
$ws.Range("D49").Value = "'16.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.08%  '
$ws.Range("D50").Value = '0.0₆0204'
$ws.Range("E50").Value = '  +2.45%  '
$ws.Range("E51").Value = '  -0.10%  '
